$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage for numeric-looking
# strings (so Excel does not silently coerce "1.00" -> 1, "0.0000174" ->
# scientific notation, etc.), then restore the default "Normal" style so
# the cell keeps no explicit style attribute, matching the source file.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "61.562.97"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.371.00"
$ws.Range("E3").Value = "  -0.89%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "572.66"
$ws.Range("E5").Value = "  -0.62%  "
Set-TextValue $ws.Range("D6") "137.23"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.371.62"
$ws.Range("E8").Value = "  -0.88%  "
Set-TextValue $ws.Range("D9") "0.471"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("E11").Value = "  -1.03%  "
Set-TextValue $ws.Range("D12") "0.392"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "3.951.24"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("E14").Value = "  +1.81%  "
Set-TextValue $ws.Range("D15") "0.0000174"
$ws.Range("E15").Value = "  -2.47%  "
Set-TextValue $ws.Range("D16") "25.94"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("D17").Value = "3.377.18"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "61.697.96"
$ws.Range("E18").Value = "  -0.09%  "
Set-TextValue $ws.Range("D19") "5.89"
$ws.Range("E19").Value = "  -0.11%  "
Set-TextValue $ws.Range("D20") "13.91"
$ws.Range("E20").Value = "  -0.66%  "
Set-TextValue $ws.Range("D21") "9.30"
$ws.Range("E21").Value = "  -1.41%  "
Set-TextValue $ws.Range("D22") "377.43"
$ws.Range("E22").Value = "  -2.78%  "
Set-TextValue $ws.Range("D23") "0.552"
$ws.Range("E23").Value = "  -3.37%  "
$ws.Range("D24").Value = "3.517.15"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("E25").Value = "  -0.11%  "
Set-TextValue $ws.Range("D26") "71.09"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  -1.38%  "
Set-TextValue $ws.Range("D28") "1.80"
$ws.Range("E28").Value = "  +8.70%  "
Set-TextValue $ws.Range("D29") "7.58"
$ws.Range("E29").Value = "  -3.13%  "
Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("E31").Value = "  +4.51%  "
Set-TextValue $ws.Range("D32") "8.20"
$ws.Range("E32").Value = "  -1.12%  "
Set-TextValue $ws.Range("D33") "2.15"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  +0.05%  "
Set-TextValue $ws.Range("D35") "23.56"
$ws.Range("E35").Value = "  +0.00%  "
Set-TextValue $ws.Range("D36") "5.20"
$ws.Range("E36").Value = "  -5.72%  "
Set-TextValue $ws.Range("D37") "1.54"
$ws.Range("E37").Value = "  -0.94%  "
Set-TextValue $ws.Range("D38") "6.81"
$ws.Range("E38").Value = "  -3.44%  "
Set-TextValue $ws.Range("D39") "164.68"
$ws.Range("E39").Value = "  +2.04%  "
Set-TextValue $ws.Range("D40") "0.0766"
$ws.Range("E40").Value = "  -3.67%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D41") "1.00"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D42") "1.73"
$ws.Range("E42").Value = "  +0.39%  "
Set-TextValue $ws.Range("D43") "0.773"
$ws.Range("E43").Value = "  +0.27%  "
Set-TextValue $ws.Range("D44") "1.21"
$ws.Range("E44").Value = "  -0.98%  "
Set-TextValue $ws.Range("D45") "41.46"
$ws.Range("E45").Value = "  +0.32%  "
Set-TextValue $ws.Range("D46") "4.36"
$ws.Range("E46").Value = "  -1.76%  "
Set-TextValue $ws.Range("D47") "24.17"
$ws.Range("E47").Value = "  +3.46%  "
Set-TextValue $ws.Range("D48") "6.84"
$ws.Range("E48").Value = "  -2.27%  "
Set-TextValue $ws.Range("D49") "22.91"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "2.373.66"
$ws.Range("E50").Value = "  +1.36%  "
Set-TextValue $ws.Range("D51") "0.0262"
$ws.Range("E51").Value = "  -1.86%  "
